# Update Name of Algo
# Apply updated values produced by the (re-)run of the RandomForest imputation
# algorithm for a handful of cells in column C (and one in column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -12.8139
$ws.Range("C12").Value = -14.16450000000001
$ws.Range("D13").Value = -8.065100000000001
$ws.Range("C18").Value = -13.9659
